$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 150 (pushes old rows 150..200 down to 151..201)
$ws.Rows.Item(150).Insert()

# Populate the newly inserted row 150 with the new record
$ws.Cells.Item(150, 1).Value = 10
$ws.Cells.Item(150, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(150, 3).Value = "La Araucanía"
$ws.Cells.Item(150, 4).Value = 45027
$ws.Cells.Item(150, 5).Value = 9
$ws.Cells.Item(150, 6).Value = "Fruta"
$ws.Cells.Item(150, 7).Value = 100104
$ws.Cells.Item(150, 8).Value = "Frutos de pepita"
$ws.Cells.Item(150, 9).Value = 100104001
$ws.Cells.Item(150, 10).Value = "Granada"
$ws.Cells.Item(150, 11).Value = "Wonderfull"
$ws.Cells.Item(150, 12).Value = "Primera"
$ws.Cells.Item(150, 13).Value = 35
$ws.Cells.Item(150, 14).Value = 25000
$ws.Cells.Item(150, 15).Value = 25000
$ws.Cells.Item(150, 16).Value = 25000
$ws.Cells.Item(150, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(150, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(150, 19).Value = 1667
$ws.Cells.Item(150, 20).Value = 15
